# Update sheet1 with newly computed TPM-based NATMI values.
# The source table (Sending cluster / Ligand / Receptor / Target cluster / metrics)
# is regenerated: it now also includes the previously-missing "-> ECs" target-cluster
# rows, in addition to the existing "-> MuSCs" rows, so the data block grows from
# 3 rows to 6 rows (plus the header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs",   "Rtn4", "Tnfrsf19", "ECs",   3, 1, 13.32779766666667,  39.983393,          0.1697233513642653, 0.1697233513642653, 3, 1, 0.1440466666666667, 0.43214,   0.1606272543719576, 0.1606272543719576, 1.919824827891111,  17.27842345102,     0.02726219593244899, 0.02726219593244899),
    @("ECs",   "Rtn4", "Tnfrsf19", "MuSCs", 3, 1, 13.32779766666667,  39.983393,          0.1697233513642653, 0.1697233513642653, 3, 1, 0.7527293333333334, 2.258188,  0.8393727456280424, 0.8393727456280424, 10.03222425243156,  90.29001827188401,  0.1424611554318163,  0.1424611554318163),
    @("FAPs",  "Rtn4", "Tnfrsf19", "ECs",   3, 1, 43.30706799999999,  129.921204,         0.5514955210569645, 0.5514955210569645, 3, 1, 0.1440466666666667, 0.43214,   0.1606272543719576, 0.1606272543719576, 6.238238788506666,  56.14414909656,     0.08858521134581236, 0.08858521134581236),
    @("FAPs",  "Rtn4", "Tnfrsf19", "MuSCs", 3, 1, 43.30706799999999,  129.921204,         0.5514955210569645, 0.5514955210569645, 3, 1, 0.7527293333333334, 2.258188,  0.8393727456280424, 0.8393727456280424, 32.59850042426133,  293.386503818352,   0.4629103097111522,  0.4629103097111522),
    @("MuSCs", "Rtn4", "Tnfrsf19", "ECs",   3, 1, 21.891734,          65.67520200000001,  0.27878112757877,   0.27878112757877,   3, 1, 0.1440466666666667, 0.43214,   0.1606272543719576, 0.1606272543719576, 3.153431310253334,  28.38088179228001,  0.04477984709369627, 0.04477984709369627),
    @("MuSCs", "Rtn4", "Tnfrsf19", "MuSCs", 3, 1, 21.891734,          65.67520200000001,  0.27878112757877,   0.27878112757877,   3, 1, 0.7527293333333334, 2.258188,  0.8393727456280424, 0.8393727456280424, 16.47855033933067,  148.306953053976,   0.2340012804850738,  0.2340012804850738)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
